# Updated script and tests
$wb = $excel.ActiveWorkbook
$wsProperty = $wb.Worksheets.Item("Property")
$wsBroker   = $wb.Worksheets.Item("Broker")

# Property!O2 (DateAgreement) moves from 1/1/2020 to 9/1/2020
$wsProperty.Range("O2").Value = 44075

# Broker sheet: last remembered selection becomes C9 (single cell)
$wsBroker.Activate() | Out-Null
$wsBroker.Range("C9").Select() | Out-Null

# Leave "Property" as the active sheet/tab, with O2 selected and scrolled into view
$wsProperty.Activate() | Out-Null
$wsProperty.Range("O2").Select() | Out-Null
